$wb = $excel.ActiveWorkbook

# "展览" (Exhibitions) sheet - first section in diff
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 335
$ws1.Range("F4").Value = 423
$ws1.Range("F5").Value = 1720
$ws1.Range("F6").Value = 84
$ws1.Range("F7").Value = 2174
$ws1.Range("F11").Value = 4882
$ws1.Range("F15").Value = 227
$ws1.Range("F17").Value = 177
$ws1.Range("F18").Value = 35
$ws1.Range("F20").Value = 121
$ws1.Range("F21").Value = 3823
$ws1.Range("F22").Value = 700
$ws1.Range("F27").Value = 114
$ws1.Range("F28").Value = 21
$ws1.Range("F30").Value = 85
$ws1.Range("F31").Value = 573
$ws1.Range("F34").Value = 919
$ws1.Range("F35").Value = 2434

# "全部类型" (All types) sheet - second section in diff (row numbers shifted by 1 for the last two entries)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 335
$ws4.Range("F4").Value = 423
$ws4.Range("F5").Value = 1720
$ws4.Range("F6").Value = 84
$ws4.Range("F7").Value = 2174
$ws4.Range("F11").Value = 4882
$ws4.Range("F15").Value = 227
$ws4.Range("F17").Value = 177
$ws4.Range("F18").Value = 35
$ws4.Range("F20").Value = 121
$ws4.Range("F21").Value = 3823
$ws4.Range("F22").Value = 700
$ws4.Range("F27").Value = 114
$ws4.Range("F28").Value = 21
$ws4.Range("F30").Value = 85
$ws4.Range("F31").Value = 573
$ws4.Range("F35").Value = 919
$ws4.Range("F36").Value = 2434
